$d = $word.ActiveDocument

# Locate the unique anchor text that ends the document's existing content
# ("...low-poly-rocks") and grow three new, clean paragraphs after it by
# inserting one paragraph mark at a time. Doing this one `^p` per
# Find/Replace call (rather than a single multi-`^p` replacement) keeps
# each newly created paragraph free of stray empty runs.
$anchor = "low-poly-rocks"

$d.Content.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "low-poly-rocks^p", 2) | Out-Null
$d.Content.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "low-poly-rocks^p", 2) | Out-Null
$d.Content.Find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, "low-poly-rocks^p", 2) | Out-Null

# After the three inserts, the trailing paragraphs are:
#   ... "low-poly-rocks" paragraph (unchanged)
#   <empty paragraph>           -> stays empty
#   <empty paragraph>           -> becomes "speedometer/fuel gauge"
#   <empty paragraph>           -> becomes the new hyperlink
#   <original trailing empty paragraph>
$count = $d.Paragraphs.Count
$textParaIndex = $count - 2
$linkParaIndex = $count - 1

$textPara = $d.Paragraphs.Item($textParaIndex)
$textPara.Range.Text = "speedometer/fuel gauge"

$linkPara = $d.Paragraphs.Item($linkParaIndex)
$d.Hyperlinks.Add($linkPara.Range, "https://www.youtube.com/watch?v=3xSYkFdQiZ0&ab_channel=CodeMonkey") | Out-Null
